$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.848.28'
$ws.Range("E2").Value = '  -4.17%  '
$ws.Range("D3").Value = '2.991.86'
$ws.Range("E3").Value = '  -4.65%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''541.83'
$ws.Range("E5").Value = '  -5.36%  '
$ws.Range("D6").Value = '''152.48'
$ws.Range("E6").Value = '  -6.83%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").Value = '''0.568'
$ws.Range("E8").Value = '  -0.58%  '
$ws.Range("D9").Value = '3.005.55'
$ws.Range("E9").Value = '  -4.68%  '
$ws.Range("E10").Value = '  -3.79%  '
$ws.Range("D11").Value = '''6.16'
$ws.Range("E11").Value = '  -6.94%  '
$ws.Range("D12").Value = '''0.369'
$ws.Range("E12").Value = '  -3.40%  '
$ws.Range("D13").Value = '3.513.38'
$ws.Range("E13").Value = '  -4.76%  '
$ws.Range("D15").Value = '61.903.19'
$ws.Range("E15").Value = '  -4.06%  '
$ws.Range("D16").Value = '''23.89'
$ws.Range("E16").Value = '  -3.89%  '
$ws.Range("D17").Value = '2.998.20'
$ws.Range("E17").Value = '  -4.87%  '
$ws.Range("D18").Value = '''0.0000147'
$ws.Range("E18").Value = '  -5.28%  '
$ws.Range("D19").Value = '''5.16'
$ws.Range("E19").Value = '  -1.38%  '
$ws.Range("D20").Value = '''12.06'
$ws.Range("E20").Value = '  -3.38%  '
$ws.Range("D21").Value = '''379.53'
$ws.Range("E21").Value = '  -8.33%  '
$ws.Range("D22").Value = '''6.76'
$ws.Range("E22").Value = '  -4.07%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '  -0.18%  '
$ws.Range("D24").Value = '''65.98'
$ws.Range("E24").Value = '  -3.90%  '
$ws.Range("D25").Value = '3.118.93'
$ws.Range("E25").Value = '  -4.56%  '
$ws.Range("D26").Value = '''0.470'
$ws.Range("E26").Value = '  -2.52%  '
$ws.Range("D27").Value = '''0.189'
$ws.Range("E27").Value = '  -2.01%  '
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").Value = '''0.997'
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0941'
$ws.Range("E29").Value = '  -7.96%  '
$ws.Range("D30").Value = '''8.24'
$ws.Range("E30").Value = '  -6.96%  '
$ws.Range("E32").Value = '  -4.18%  '
$ws.Range("D33").Value = '''20.46'
$ws.Range("E33").Value = '  -3.58%  '
$ws.Range("D34").Value = '''161.19'
$ws.Range("E34").Value = '  -0.87%  '
$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").Value = '''5.98'
$ws.Range("E35").Value = '  -4.23%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '''4.61'
$ws.Range("E36").Value = '  -5.12%  '
$ws.Range("D37").Value = '''1.07'
$ws.Range("E37").Value = '  -4.72%  '
$ws.Range("D38").Value = '''1.28'
$ws.Range("E38").Value = '  -5.06%  '
$ws.Range("D39").Value = '''1.55'
$ws.Range("E39").Value = '  -7.77%  '
$ws.Range("D40").Value = '''37.70'
$ws.Range("E40").Value = '  -1.40%  '
$ws.Range("D41").Value = '2.423.95'
$ws.Range("E41").Value = '  -7.84%  '
$ws.Range("D42").Value = '''3.91'
$ws.Range("E42").Value = '  -4.87%  '
$ws.Range("D43").Value = '''22.15'
$ws.Range("E43").Value = '  -6.71%  '
$ws.Range("D44").Value = '''0.673'
$ws.Range("E44").Value = '  -2.42%  '
$ws.Range("D45").Value = '''0.0592'
$ws.Range("E45").Value = '  -3.21%  '
$ws.Range("D46").Value = '''5.17'
$ws.Range("E46").Value = '  -2.24%  '
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("D48").Value = '''0.0245'
$ws.Range("E48").Value = '  -3.47%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '''0.0954'
$ws.Range("E49").Value = '  -2.08%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '''19.83'
$ws.Range("E50").Value = '  -6.76%  '
$ws.Range("D51").Value = '''268.21'
$ws.Range("E51").Value = '  -7.45%  '
